# Adds three new worksheets (validPythonCode, submitPythonCode,
# invalidPythonCode) with "practice question" data, matching the commit
# "Added excel sheet for practice questions".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Reusable code snippets / text pulled from the existing pythonCode sheet
# so the shared strings get reused (same text => same shared string).
# ---------------------------------------------------------------------
$searchCode = "def search(input_list, num):`nif(num in input_list):`nprint(`"Element Found`")`n\b`n\b`nelse:`nprint(`"Not Found`")`n\b`n\b`n\b`n\b`nsearch([12, 23, 45, 67, 6, 90] , 12)"
$maxConsecCode = "def findMaxConsecutiveOnes(nums) :`ncount = 0`nresult = 0`nfor i in range(0, len(nums)):`nif (nums[i] == 0):`ncount = 0`n\b`n\b`nelse:`ncount+= 1`n\b`n\b`nresult = max(result, count)`n\b`n\b`nprint(result)`n\b`n\b`nfindMaxConsecutiveOnes([1,0,1,1,0,1])"
$findNumbersCode = "def findNumbers(nums):`nc=0`nfor i in nums:`nj=str(i)`nx=len(j)`nif x%2==0:`nc=c+1`n\b`n\b`n\b`n\b`nprint c`nreturn c`nfindNumbers([12,345,2,6,7896])"
$sortedSquaresCode = "def sortedSquares(nums):`nsquares_list = []`nfor i in range(0, len(nums)):`nsquare = nums[i] * nums[i];`nsquares_list.append(square)`n\b`n\b`nsorted_squares_list = sorted(squares_list)`nprint sorted_squares_list;`nreturn sorted_squares_list;`nsortedSquares([-7,-3,2,3,11])"

$elementFound = "Element Found"
$two = "2"
$submissionSuccess = "submission success"
$squaresResult = "[4, 9, 9, 49, 121]"
$hello = "hello"
$nameError = "NameError: name 'hello' is not defined on line 1"

$questionTitle = "questionTitle"
$pythonCodeHeader = "pythonCode"
$resultHeader = "Result"
$q1 = "Search the array"
$q2 = "Max Consecutive Ones"
$q3 = "Find Numbers with Even Number of Digits"
$q4 = "Squares of a Sorted Array"

# ---------------------------------------------------------------------
# Sheet: validPythonCode  (added after the last existing sheet so the
# new sheets land at the end of the tab strip, in commit order)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsValid = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsValid.Name = "validPythonCode"

$wsValid.Range("A1").Value = $questionTitle
$wsValid.Range("B1").Value = $pythonCodeHeader
$wsValid.Range("C1").Value = $resultHeader

$wsValid.Range("A2").Value = $q1
$wsValid.Range("B2").Value = $searchCode
$wsValid.Range("C2").Value = $elementFound

$wsValid.Range("A3").Value = $q2
$wsValid.Range("B3").Value = $maxConsecCode
$wsValid.Range("C3").Value = $two

$wsValid.Range("A4").Value = $q3
$wsValid.Range("B4").Value = $findNumbersCode
$wsValid.Range("C4").Value = $two

$wsValid.Range("A5").Value = $q4
$wsValid.Range("B5").Value = $sortedSquaresCode
$wsValid.Range("C5").Value = $squaresResult

$wsValid.Columns.Item(1).ColumnWidth = 39
$wsValid.Columns.Item(2).ColumnWidth = 103.28515625

# ---------------------------------------------------------------------
# Sheet: submitPythonCode
# ---------------------------------------------------------------------
$wsSubmit = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsValid)
$wsSubmit.Name = "submitPythonCode"

$wsSubmit.Range("A1").Value = $questionTitle
$wsSubmit.Range("B1").Value = $pythonCodeHeader
$wsSubmit.Range("C1").Value = $resultHeader

$wsSubmit.Range("A2").Value = $q1
$wsSubmit.Range("B2").Value = $searchCode
$wsSubmit.Range("C2").Value = $submissionSuccess

$wsSubmit.Range("A3").Value = $q2
$wsSubmit.Range("B3").Value = $maxConsecCode
$wsSubmit.Range("C3").Value = $submissionSuccess

$wsSubmit.Range("A4").Value = $q3
$wsSubmit.Range("B4").Value = $findNumbersCode
$wsSubmit.Range("C4").Value = $submissionSuccess

$wsSubmit.Range("A5").Value = $q4
$wsSubmit.Range("B5").Value = $sortedSquaresCode
$wsSubmit.Range("C5").Value = $submissionSuccess

$wsSubmit.Columns.Item(1).ColumnWidth = 39
$wsSubmit.Columns.Item(2).ColumnWidth = 118.85546875

# ---------------------------------------------------------------------
# Sheet: invalidPythonCode
# ---------------------------------------------------------------------
$wsInvalid = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSubmit)
$wsInvalid.Name = "invalidPythonCode"

$wsInvalid.Range("A1").Value = $questionTitle
$wsInvalid.Range("B1").Value = $pythonCodeHeader
$wsInvalid.Range("C1").Value = $resultHeader

$wsInvalid.Range("A2").Value = $q1
$wsInvalid.Range("B2").Value = $hello
$wsInvalid.Range("C2").Value = $nameError

$wsInvalid.Range("A3").Value = $q2
$wsInvalid.Range("B3").Value = $hello
$wsInvalid.Range("C3").Value = $nameError

$wsInvalid.Range("A4").Value = $q3
$wsInvalid.Range("B4").Value = $hello
$wsInvalid.Range("C4").Value = $nameError

$wsInvalid.Range("A5").Value = $q4
$wsInvalid.Range("B5").Value = $hello
$wsInvalid.Range("C5").Value = $nameError

$wsInvalid.Columns.Item(1).ColumnWidth = 39
$wsInvalid.Columns.Item(2).ColumnWidth = 16.140625
$wsInvalid.Columns.Item(3).ColumnWidth = 44.28515625

$wsInvalid.Activate()
